$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2- crescimento")

# Rename "Não" -> "Nao" in the existing "Desmamada" answers (rows 2-16),
# and rename the "Data nasc." header to "Data Nasc." — order matters so the
# shared-string table is rebuilt with "Sim","Nao","Data Nasc." appended in
# that sequence.
$ws.Range("F2:F16").Value = "Nao"
$ws.Range("B1").Value = "Data Nasc."

# Extend the growth tracking table with 10 more blank rows (19-28), giving
# column F (Desmamada) the same number formatting/style as the rows above.
$ws.Range("F19:F28").NumberFormat = $ws.Range("F2").NumberFormat

# Add the dropdown list validations for the "Desmamada" column.
$ws.Range("F2:F18").Validation.Add(3, 1, 1, '"Sim, Nao"')
$ws.Range("F28").Validation.Add(3, 1, 1, '"Sim, Nao, Vazio"')

# Add a threaded comment on the "Idade ao Desm." header explaining usage.
$excel.UserName = "Tadeu Da Silva"
$ws.Range("G1").AddCommentThreaded("Caso a bezerra não tenha sido desaleitada, deixe vazio, ou seja, não adicione nenhum texto à célula.") | Out-Null

# Restore the frozen-pane view with the new active selection.
$ws.Activate()
$ws.Range("B2").Select() | Out-Null
